$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$rng = $ws1.Range("Z1")
$rng.Interior.Pattern = 1
$rng.Interior.Color = 33023
Write-Host (Get-Member -InputObject $rng.Interior)
